# Generate Report for Handoff
# Updates the "Status" text from "In Translation" to "Ready for handoff"
# and refreshes the handoff timestamps on all three sheets, then widens
# the columns that display the status text so the new wording fits.

$wb = $excel.ActiveWorkbook

# Target column width (OOXML "width" attribute) is 17.2159881591797 chars,
# widened from 13.4101845877511 to fit the new "Ready for handoff" text.
# The host quantizes ColumnWidth writes to the nearest 1/6 of a character,
# so 16.3 is the closest settable value that lands on the same bucket.
$newColWidth = 16.3

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-25 02:58:09"
$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-25 02:57:57"
$wsZhCn.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-25 02:58:09"
$wsDeDe.Columns.Item(3).ColumnWidth = $newColWidth
